# Weekly update: prepend 4 new rows of data (new week) to the Tomate sheet,
# pushing the existing data (rows 1140:1226) down to (1144:1230).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current row 1140, shifting everything
# from old row 1140 onward down by 4 rows.
$ws.Rows("1140:1143").Insert()

# Fill in the 4 new rows with this week's data.

# Row 1140
$ws.Range("A1140").Value = 10
$ws.Range("B1140").Value = "Vega Modelo de Temuco"
$ws.Range("C1140").Value = "La Araucanía"
$ws.Range("D1140").Value = 44585
$ws.Range("E1140").Value = 9
$ws.Range("F1140").Value = 100112020
$ws.Range("G1140").Value = "Tomate"
$ws.Range("H1140").Value = "Larga vida"
$ws.Range("I1140").Value = "Extra"
$ws.Range("J1140").Value = 550
$ws.Range("K1140").Value = 18000
$ws.Range("L1140").Value = 18000
$ws.Range("M1140").Value = 18000
$ws.Range("N1140").Value = "`$/bandeja 18 kilos"
$ws.Range("O1140").Value = "Limache"
$ws.Range("P1140").Value = 1000
$ws.Range("Q1140").Value = 18
$ws.Range("R1140").Value = "Hortaliza"

# Row 1141
$ws.Range("A1141").Value = 10
$ws.Range("B1141").Value = "Vega Modelo de Temuco"
$ws.Range("C1141").Value = "La Araucanía"
$ws.Range("D1141").Value = 44585
$ws.Range("E1141").Value = 9
$ws.Range("F1141").Value = 100112020
$ws.Range("G1141").Value = "Tomate"
$ws.Range("H1141").Value = "Larga vida"
$ws.Range("I1141").Value = "Primera"
$ws.Range("J1141").Value = 650
$ws.Range("K1141").Value = 12000
$ws.Range("L1141").Value = 13000
$ws.Range("M1141").Value = 12462
$ws.Range("N1141").Value = "`$/bandeja 18 kilos"
$ws.Range("O1141").Value = "Angol"
$ws.Range("P1141").Value = 692
$ws.Range("Q1141").Value = 18
$ws.Range("R1141").Value = "Hortaliza"

# Row 1142
$ws.Range("A1142").Value = 10
$ws.Range("B1142").Value = "Vega Modelo de Temuco"
$ws.Range("C1142").Value = "La Araucanía"
$ws.Range("D1142").Value = 44585
$ws.Range("E1142").Value = 9
$ws.Range("F1142").Value = 100112020
$ws.Range("G1142").Value = "Tomate"
$ws.Range("H1142").Value = "Larga vida"
$ws.Range("I1142").Value = "Primera"
$ws.Range("J1142").Value = 550
$ws.Range("K1142").Value = 14000
$ws.Range("L1142").Value = 15000
$ws.Range("M1142").Value = 14455
$ws.Range("N1142").Value = "`$/bandeja 18 kilos"
$ws.Range("O1142").Value = "Limache"
$ws.Range("P1142").Value = 803
$ws.Range("Q1142").Value = 18
$ws.Range("R1142").Value = "Hortaliza"

# Row 1143
$ws.Range("A1143").Value = 10
$ws.Range("B1143").Value = "Vega Modelo de Temuco"
$ws.Range("C1143").Value = "La Araucanía"
$ws.Range("D1143").Value = 44585
$ws.Range("E1143").Value = 9
$ws.Range("F1143").Value = 100112020
$ws.Range("G1143").Value = "Tomate"
$ws.Range("H1143").Value = "Semiduro"
$ws.Range("I1143").Value = "Primera"
$ws.Range("J1143").Value = 600
$ws.Range("K1143").Value = 10000
$ws.Range("L1143").Value = 11000
$ws.Range("M1143").Value = 10583
$ws.Range("N1143").Value = "`$/bandeja 18 kilos"
$ws.Range("O1143").Value = "Región del Maule"
$ws.Range("P1143").Value = 588
$ws.Range("Q1143").Value = 18
$ws.Range("R1143").Value = "Hortaliza"
